$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name shown in workbook.xml <sheet name=...>)
$ws.Name = "BetaFiberA"

# Tiny floating point precision corrections (last-bit rounding) in row 13 and row 15
$ws.Range("C13").Value = 0.9957355901491003
$ws.Range("G13").Value = 0.9957355901491003
$ws.Range("M13").Value = 0.9957355901491003
$ws.Range("O13").Value = 0.9948144489558475
$ws.Range("K15").Value = 0.9850765749113386

# New row 16 of data (HKL index 14, label "HexGrid-60degTilt5degRes")
# Copy formatting from the row above (A15/B15) so A16 keeps the same bold/border/center style.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.062027823678512
$ws.Range("D16").Value = 0.9128888380768491
$ws.Range("E16").Value = 1.011713906380467
$ws.Range("F16").Value = 0.973209509643741
$ws.Range("G16").Value = 1.062027823678512
$ws.Range("H16").Value = 0.9128888380768491
$ws.Range("I16").Value = 1.02520033075858
$ws.Range("J16").Value = 0.9715609114691947
$ws.Range("K16").Value = 1.015461794302678
$ws.Range("L16").Value = 0.9421122018797425
$ws.Range("M16").Value = 1.062027823678512
$ws.Range("N16").Value = 0.9623013722286582
$ws.Range("O16").Value = 0.9899600194448923
$ws.Range("P16").Value = 0.9892719145237204
